# Apply the "PO Forecast" sheet addition plus header-text fixes.
$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# 1) Rename the "Requested quantity" headers.
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add a new worksheet "PO Forecast" after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the bold/bordered header style used on the other sheets.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Copy the date-number-format style used in column A of the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122)  # xlPasteFormats

$rows = @(
    @(2, 45550.99999999999, 211, -240.8749976699012, 700.9558686550583),
    @(3, 45578.99999999999, 546, 89.22842312848988, 1028.285497237619),
    @(4, 45592.99999999999, 714, 278.1852029953668, 1140.688965919013),
    @(5, 45599.99999999999, 798, 359.6586675257141, 1244.510221738313),
    @(6, 45620.99999999999, 1050, 594.4651964612233, 1511.124196653266),
    @(7, 45627.99999999999, 1134, 694.3326761466967, 1604.243709392822),
    @(8, 45634.99999999999, 1218, 718.4781133624289, 1695.580167203426),
    @(9, 45641.99999999999, 1302, 877.1615635128177, 1763.649341992555),
    @(10, 45648.99999999999, 1386, 888.6859781212025, 1821.667934475106),
    @(11, 45655.99999999999, 1470, 982.8833799595752, 1924.194531438166),
    @(12, 45662.99999999999, 1554, 1080.224910078755, 1998.333049952303),
    @(13, 45669.99999999999, 1638, 1220.154922523033, 2120.221252926188),
    @(14, 45676.99999999999, 1722, 1266.118946085651, 2169.805987681213)
)

foreach ($r in $rows) {
    $row = $r[0]
    $wsForecast.Cells.Item($row, 1).Value = $r[1]
    $wsForecast.Cells.Item($row, 2).Value = $r[2]
    $wsForecast.Cells.Item($row, 3).Value = $r[3]
    $wsForecast.Cells.Item($row, 4).Value = $r[4]
}

# Restore the originally-active sheet/selection (adding a sheet makes it active).
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
